# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) was previously populated from a "Strike#" (total
# strikes thrown) source; the pipeline was regenerated to populate it from
# strikeouts (K) instead, so every row's G value is rewritten here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) counts for column G, rows 2-36, in row order.
$sVals = @(3, 6, 2, 6, 6, 5, 2, 1, 3, 6, 6, 7, 6, 3, 2, 6, 8, 7, 3, 1, 6, 2, 4, 3, 3, 3, 8, 4, 10, 4, 6, 6, 3, 1, 1)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
